# Updated cryptos list on Fri May 19 15:00:30 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.847.60"
$ws.Range('E2').Value = "'  -1.48%  "
$ws.Range('D3').Value = "'1.811.69"
$ws.Range('E3').Value = "'  -0.36%  "
$ws.Range('D4').Value = "'1.001"
$ws.Range('E4').Value = "'  -0.01%  "
$ws.Range('D5').Value = "'310.29"
$ws.Range('E5').Value = "'  -0.86%  "
$ws.Range('E6').Value = "'  -0.01%  "
$ws.Range('D7').Value = "'0.4710"
$ws.Range('E7').Value = "'  +1.53%  "
$ws.Range('D8').Value = "'0.3699"
$ws.Range('E8').Value = "'  -1.54%  "
$ws.Range('D9').Value = "'0.07351"
$ws.Range('E9').Value = "'  -0.82%  "
$ws.Range('D10').Value = "'0.8687"
$ws.Range('E10').Value = "'  +0.04%  "
$ws.Range('D11').Value = "'20.39"
$ws.Range('E11').Value = "'  -1.09%  "
$ws.Range('D12').Value = "'1.910.01"
$ws.Range('E12').Value = "'  +5.00%  "
$ws.Range('D13').Value = "'5.348"
$ws.Range('E13').Value = "'  -0.97%  "
$ws.Range('E14').Value = "'  -0.25%  "
$ws.Range('D15').Value = "'6.503"
$ws.Range('E15').Value = "'  -2.43%  "
$ws.Range('D16').Value = "'91.66"
$ws.Range('E16').Value = "'  -0.48%  "
$ws.Range('D17').Value = "'1.002"
$ws.Range('E17').Value = "'  +0.00%  "
$ws.Range('D18').Value = "'0.000008698"
$ws.Range('E18').Value = "'  -0.73%  "
$ws.Range('D19').Value = "'1.001"
$ws.Range('E19').Value = "'  +0.03%  "
$ws.Range('E20').Value = "'  -1.65%  "
$ws.Range('D21').Value = "'26.887.45"
$ws.Range('D22').Value = "'5.333"
$ws.Range('E22').Value = "'  +0.35%  "
$ws.Range('D23').Value = "'10.54"
$ws.Range('E23').Value = "'  -3.42%  "
$ws.Range('D24').Value = "'2.099.95"
$ws.Range('E24').Value = "'  +2.48%  "
$ws.Range('D25').Value = "'1.895"
$ws.Range('E25').Value = "'  -2.00%  "
$ws.Range('D26').Value = "'151.88"
$ws.Range('E26').Value = "'  -0.03%  "
$ws.Range('E27').Value = "'  -0.63%  "
$ws.Range('D28').Value = "'2.097"
$ws.Range('E28').Value = "'  -7.50%  "
$ws.Range('D29').Value = "'5.288"
$ws.Range('E29').Value = "'  -0.02%  "
$ws.Range('D30').Value = "'115.28"
$ws.Range('E30').Value = "'  -1.61%  "
$ws.Range('D31').Value = "'0.08949"
$ws.Range('E31').Value = "'  +0.60%  "
$ws.Range('D32').Value = "'0.7557"
$ws.Range('E32').Value = "'  -2.82%  "
$ws.Range('B33').Value = "'ARBITRUM"
$ws.Range('C33').Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range('D33').Value = "'1.149"
$ws.Range('E33').Value = "'  -2.75%  "
$ws.Range('B34').Value = "'HuobiToken"
$ws.Range('C34').Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range('D34').Value = "'2.927"
$ws.Range('E34').Value = "'  +0.49%  "
$ws.Range('E36').Value = "'  -0.02%  "
$ws.Range('D37').Value = "'1.094"
$ws.Range('E37').Value = "'  -0.96%  "
$ws.Range('D38').Value = "'0.01952"
$ws.Range('E38').Value = "'  -0.52%  "
$ws.Range('E39').Value = "'  +0.16%  "
$ws.Range('D40').Value = "'0.5328"
$ws.Range('E40').Value = "'  +0.76%  "
$ws.Range('D41').Value = "'2.909"
$ws.Range('E41').Value = "'  +0.27%  "
$ws.Range('D42').Value = "'7.168"
$ws.Range('E42').Value = "'  -1.30%  "
$ws.Range('D43').Value = "'2.350"
$ws.Range('E43').Value = "'  -1.42%  "
$ws.Range('D44').Value = "'0.1658"
$ws.Range('E44').Value = "'  -1.84%  "
$ws.Range('D45').Value = "'8.420"
$ws.Range('E45').Value = "'  -2.36%  "
$ws.Range('D46').Value = "'0.4931"
$ws.Range('E46').Value = "'  -2.16%  "
$ws.Range('D47').Value = "'10.30"
$ws.Range('E47').Value = "'  -1.78%  "
$ws.Range('E48').Value = "'  +0.02%  "
$ws.Range('D49').Value = "'1.672"
$ws.Range('E49').Value = "'  -0.08%  "
$ws.Range('D50').Value = "'103.04"
$ws.Range('E50').Value = "'  -2.24%  "
$ws.Range('D51').Value = "'0.06275"
$ws.Range('E51').Value = "'  -0.74%  "
